$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Adjust the "subno" column (B) from 696 to 682 for every data row (2-54)
# to match the folder number (s682_1 / DD_682).
$ws.Range("B2:B54").Value = 682

# Restore the window/selection state recorded after the edit: scroll so
# row 24 is at the top of the view, with B2:B54 selected (active cell B2).
$win = $wb.Windows.Item(1)
[void]$ws.Range("B2:B54").Select()
$win.ScrollRow = 24
$win.ScrollColumn = 1
